# Flip the sign of every non-zero numeric value in the B2:E5 data block
# ("material recycled" per component) on every year sheet of the workbook.
$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    for ($r = 2; $r -le 5; $r++) {
        for ($c = 2; $c -le 5; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne $null -and $val -ne 0) {
                $cell.Value = -$val
            }
        }
    }
}
